# The pre-seeded $wb variable does not support method calls in this host
# (property reads work, but `$wb.Foo()` / chained `.Method()` calls throw
# "You cannot call a method on a null-valued expression."). Re-fetching the
# active workbook through $excel works fine for both properties and methods.
$wb2 = $excel.ActiveWorkbook
$ws  = $wb2.ActiveSheet

# delete_other_columns: drop a column outside the real A:E data block
# (mirrors the commit's columns_from_strings-driven column deletion; the
# worksheet's used range A1:E10 is untouched, only the far-right custom
# column-width bands shift left by one as a side effect).
$ws.Range("F1").EntireColumn.Delete()

# Rename the image filenames referenced in column C.
$ws.Range("C2").Value = "congruent1_o.jpg"
$ws.Range("C3").Value = "incongruent1_n.jpg"
$ws.Range("C4").Value = "congruent2_o.jpg"
$ws.Range("C5").Value = "congruent3_n.jpg"
$ws.Range("C6").Value = "congruent5_o.jpg"
$ws.Range("C7").Value = "d1_o.jpg"
$ws.Range("C8").Value = "congruent1_o.jpg"
$ws.Range("C9").Value = "neutre12_n.jpg"
